$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("H1").Value = "Consignee Name"
$ws.Range("I1").Value = "Destination"
$ws.Range("J1").Value = "Dest Branch"
$ws.Range("K1").Value = "Product"
$ws.Range("L1").Value = "Quality"
$ws.Range("M1").Value = "D1"
$ws.Range("N1").Value = "D2"
$ws.Range("O1").Value = "D4"
$ws.Range("P1").Value = "Booking Type"
$ws.Range("Q1").Value = "Customer Type"
$ws.Range("R1").Value = "Region"
$ws.Range("S1").Value = "Contract Price"
$ws.Range("T1").Value = "Planned Quantity"
$ws.Range("U1").Value = "TDCs"
$ws.Range("V1").Value = "Delivery Date"
$ws.Range("W1").Value = "Rail/Road Ind"

# Row 2
$ws.Range("H2").Value = "D"
$ws.Range("I2").Value = "I"
$ws.Range("J2").Value = "B034(BSO BARODA)"
$ws.Range("K2").Value = "HR COILS"
$ws.Range("L2").Value = "IS 2062 E250BR"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = "19000010"
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = "MOU"
$ws.Range("Q2").Value = "GEN"
$ws.Range("R2").Value = "Western Region"
$ws.Range("S2").Value = "39200"
$ws.Range("T2").Value = "64"
$ws.Range("U2").Value = "NS001:D_E250BR"
$ws.Range("V2").Value = "31.01.2025"
$ws.Range("W2").Value = "DD - Rail"

# Row 3
$ws.Range("H3").Value = "D"
$ws.Range("I3").Value = "I"
$ws.Range("J3").Value = "B034(BSO BARODA)"
$ws.Range("K3").Value = "HR COILS"
$ws.Range("L3").Value = "IS 2062 E250BR"
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = "19000010"
$ws.Range("O3").Value = ""
$ws.Range("P3").Value = "MOU"
$ws.Range("Q3").Value = "GEN"
$ws.Range("R3").Value = "Western Region"
$ws.Range("S3").Value = "39200"
$ws.Range("T3").Value = "64"
$ws.Range("U3").Value = "NS001:D_E250BR"
$ws.Range("V3").Value = "31.01.2025"
$ws.Range("W3").Value = "DD - Rail"

# Row 4
$ws.Range("H4").Value = "D"
$ws.Range("I4").Value = "I"
$ws.Range("J4").Value = "B034(BSO BARODA)"
$ws.Range("K4").Value = "HR COILS"
$ws.Range("L4").Value = "IS 2062 E250BR"
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = "19000010"
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = "MOU"
$ws.Range("Q4").Value = "GEN"
$ws.Range("R4").Value = "Western Region"
$ws.Range("S4").Value = "39200"
$ws.Range("T4").Value = "64"
$ws.Range("U4").Value = "NS001:D_E250BR"
$ws.Range("V4").Value = "31.01.2025"
$ws.Range("W4").Value = "DD - Rail"

# Row 5
$ws.Range("H5").Value = "D"
$ws.Range("I5").Value = "I"
$ws.Range("J5").Value = "B033(BSO AHMEDABAD)"
$ws.Range("K5").Value = "HR COILS"
$ws.Range("L5").Value = "IS 2062 E350 BR"
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = "19000015"
$ws.Range("O5").Value = ""
$ws.Range("P5").Value = "BRQ"
$ws.Range("Q5").Value = "BRM"
$ws.Range("R5").Value = "Western Region"
$ws.Range("S5").Value = "41870"
$ws.Range("T5").Value = "128"
$ws.Range("U5").Value = "NS001:D_E250BR"
$ws.Range("V5").Value = "31.01.2025"
$ws.Range("W5").Value = "Rail - WH"

# Row 6
$ws.Range("H6").Value = "D"
$ws.Range("I6").Value = "I"
$ws.Range("J6").Value = "B033(BSO AHMEDABAD)"
$ws.Range("K6").Value = "HR COILS"
$ws.Range("L6").Value = "IS 2062 E350 BR"
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = "19000015"
$ws.Range("O6").Value = ""
$ws.Range("P6").Value = "BRQ"
$ws.Range("Q6").Value = "BRM"
$ws.Range("R6").Value = "Western Region"
$ws.Range("S6").Value = "41870"
$ws.Range("T6").Value = "128"
$ws.Range("U6").Value = "NS001:D_E250BR"
$ws.Range("V6").Value = "31.01.2025"
$ws.Range("W6").Value = "Rail - WH"
